$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: sort the 6/13 block (rows 356:370) by Staff Name (column A) ascending ---
$srt = $ws.Sort
$srt.SortFields.Clear()
$srt.SortFields.Add($ws.Range("A356:A370"))
$srt.SetRange($ws.Range("A356:F370"))
$srt.Header = 0
$srt.Apply()

# --- Step 2: the sort pushes every "AV Shutdown" row (now rows 356:363) to the top; ---
# --- those shutdown tasks are done with, so remove that whole block ---
$ws.Range("A356:F363").EntireRow.Delete()

# --- Step 3: log the next day's (6/14) tasks as new rows at the bottom of the log ---
$newRows = @(
    @("Pickup Mic", 42614, "1630", "YL", "280N", "Return mic (IR) to KT 516 and place battery in charger"),
    @("AV Shutdown", 42614, "1630", "YL", "280N", ""),
    @("Demo", 42614, "1630", "OSG", "1001", ""),
    @("Demo", 42614, "1630", "OSG", "2001", ""),
    @("Demo", 42614, "1630", "OSG", "2002", ""),
    @("Demo", 42614, "1630", "OSG", "2003", ""),
    @("Demo", 42614, "1830", "OSG", "1008", "")
)

$r = 366
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    if ($row[5] -ne "") {
        $ws.Cells.Item($r, 6).Value = $row[5]
    }
    $r = $r + 1
}

# --- Step 4: leave the view/selection the way the author left it ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 357
$ws.Range("A367:XFD367").Select()
